$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.628848666666667
$ws.Range("H2").Value = 7.886546000000001
$ws.Range("I2").Value = 0.04622248078033103
$ws.Range("J2").Value = 0.04850184447997802
$ws.Range("M2").Value = 36.48539666666667
$ws.Range("N2").Value = 109.45619
$ws.Range("O2").Value = 0.4260639713374229
$ws.Range("P2").Value = 0.4324607845540777
$ws.Range("Q2").Value = 95.91458637997111
$ws.Range("R2").Value = 863.2312774197402
$ws.Range("S2").Value = 0.01969373372633554
$ws.Range("T2").Value = 0.02097514571613116
$ws.Range("G3").Value = 2.628848666666667
$ws.Range("H3").Value = 7.886546000000001
$ws.Range("I3").Value = 0.04622248078033103
$ws.Range("J3").Value = 0.04850184447997802
$ws.Range("O3").Value = 0.1743777127077069
$ws.Range("P3").Value = 0.1769957741547643
$ws.Range("Q3").Value = 39.2555280742089
$ws.Range("R3").Value = 353.2997526678801
$ws.Range("S3").Value = 0.008060170474150067
$ws.Range("T3").Value = 0.008584621511667693
$ws.Range("G4").Value = 2.628848666666667
$ws.Range("H4").Value = 7.886546000000001
$ws.Range("I4").Value = 0.04622248078033103
$ws.Range("J4").Value = 0.04850184447997802
$ws.Range("M4").Value = 9.680823666666667
$ws.Range("N4").Value = 29.042471
$ws.Range("O4").Value = 0.1130493445068016
$ws.Range("P4").Value = 0.1147466378470605
$ws.Range("Q4").Value = 25.44942038835178
$ws.Range("R4").Value = 229.044783495166
$ws.Range("S4").Value = 0.005225421153694659
$ws.Range("T4").Value = 0.005565423583458491
$ws.Range("G5").Value = 2.628848666666667
$ws.Range("H5").Value = 7.886546000000001
$ws.Range("I5").Value = 0.04622248078033103
$ws.Range("J5").Value = 0.04850184447997802
$ws.Range("M5").Value = 3.79999
$ws.Range("N5").Value = 7.59998
$ws.Range("O5").Value = 0.04437498227672168
$ws.Range("P5").Value = 0.0300274777826206
$ws.Range("Q5").Value = 9.989598644846668
$ws.Range("R5").Value = 59.93759186908001
$ws.Range("S5").Value = 0.002051121765413298
$ws.Range("T5").Value = 0.00145638805753866
$ws.Range("G6").Value = 2.628848666666667
$ws.Range("H6").Value = 7.886546000000001
$ws.Range("I6").Value = 0.04622248078033103
$ws.Range("J6").Value = 0.04850184447997802
$ws.Range("M6").Value = 20.734808
$ws.Range("N6").Value = 62.204424
$ws.Range("O6").Value = 0.242133989171347
$ws.Range("P6").Value = 0.245769325661477
$ws.Range("Q6").Value = 54.50867236438934
$ws.Range("R6").Value = 490.5780512795041
$ws.Range("S6").Value = 0.01119203366073747
$ws.Range("T6").Value = 0.01192026561118203
$ws.Range("I7").Value = 0.472133375270229
$ws.Range("J7").Value = 0.4954156322762335
$ws.Range("M7").Value = 36.48539666666667
$ws.Range("N7").Value = 109.45619
$ws.Range("O7").Value = 0.4260639713374229
$ws.Range("P7").Value = 0.4324607845540777
$ws.Range("Q7").Value = 979.706771266451
$ws.Range("R7").Value = 8817.360941398059
$ws.Range("S7").Value = 0.2011590208685756
$ws.Range("T7").Value = 0.2142478330145344
$ws.Range("I8").Value = 0.472133375270229
$ws.Range("J8").Value = 0.4954156322762335
$ws.Range("O8").Value = 0.1743777127077069
$ws.Range("P8").Value = 0.1769957741547643
$ws.Range("S8").Value = 0.08232953807259194
$ws.Range("T8").Value = 0.08768647336310399
$ws.Range("I9").Value = 0.472133375270229
$ws.Range("J9").Value = 0.4954156322762335
$ws.Range("M9").Value = 9.680823666666667
$ws.Range("N9").Value = 29.042471
$ws.Range("O9").Value = 0.1130493445068016
$ws.Range("P9").Value = 0.1147466378470605
$ws.Range("Q9").Value = 259.9497158909838
$ws.Range("R9").Value = 2339.547443018854
$ws.Range("S9").Value = 0.05337436859408318
$ws.Range("T9").Value = 0.05684727814057348
$ws.Range("I10").Value = 0.472133375270229
$ws.Range("J10").Value = 0.4954156322762335
$ws.Range("M10").Value = 3.79999
$ws.Range("N10").Value = 7.59998
$ws.Range("O10").Value = 0.04437498227672168
$ws.Range("P10").Value = 0.0300274777826206
$ws.Range("Q10").Value = 102.0374252130867
$ws.Range("R10").Value = 612.22455127852
$ws.Range("S10").Value = 0.0209509101598652
$ws.Range("T10").Value = 0.01487608189133754
$ws.Range("I11").Value = 0.472133375270229
$ws.Range("J11").Value = 0.4954156322762335
$ws.Range("M11").Value = 20.734808
$ws.Range("N11").Value = 62.204424
$ws.Range("O11").Value = 0.242133989171347
$ws.Range("P11").Value = 0.245769325661477
$ws.Range("Q11").Value = 556.7715758745974
$ws.Range("R11").Value = 5010.944182871376
$ws.Range("S11").Value = 0.1143195375751132
$ws.Range("T11").Value = 0.1217579658666842
$ws.Range("G12").Value = 11.96574466666667
$ws.Range("H12").Value = 35.897234
$ws.Range("I12").Value = 0.2103911152781009
$ws.Range("J12").Value = 0.2207661073338543
$ws.Range("M12").Value = 36.48539666666667
$ws.Range("N12").Value = 109.45619
$ws.Range("O12").Value = 0.4260639713374229
$ws.Range("P12").Value = 0.4324607845540777
$ws.Range("Q12").Value = 436.5749405753844
$ws.Range("R12").Value = 3929.17446517846
$ws.Range("S12").Value = 0.0896400741094972
$ws.Range("T12").Value = 0.09547268398054834
$ws.Range("G13").Value = 11.96574466666667
$ws.Range("H13").Value = 35.897234
$ws.Range("I13").Value = 0.2103911152781009
$ws.Range("J13").Value = 0.2207661073338543
$ws.Range("O13").Value = 0.1743777127077069
$ws.Range("P13").Value = 0.1769957741547643
$ws.Range("Q13").Value = 178.6795990378356
$ws.Range("R13").Value = 1608.11639134052
$ws.Range("S13").Value = 0.03668752145621871
$ws.Range("T13").Value = 0.03907466807468933
$ws.Range("G14").Value = 11.96574466666667
$ws.Range("H14").Value = 35.897234
$ws.Range("I14").Value = 0.2103911152781009
$ws.Range("J14").Value = 0.2207661073338543
$ws.Range("M14").Value = 9.680823666666667
$ws.Range("N14").Value = 29.042471
$ws.Range("O14").Value = 0.1130493445068016
$ws.Range("P14").Value = 0.1147466378470605
$ws.Range("Q14").Value = 115.8382641583571
$ws.Range("R14").Value = 1042.544377425214
$ws.Range("S14").Value = 0.02378457767224424
$ws.Range("T14").Value = 0.02533216856714307
$ws.Range("G15").Value = 11.96574466666667
$ws.Range("H15").Value = 35.897234
$ws.Range("I15").Value = 0.2103911152781009
$ws.Range("J15").Value = 0.2207661073338543
$ws.Range("M15").Value = 3.79999
$ws.Range("N15").Value = 7.59998
$ws.Range("O15").Value = 0.04437498227672168
$ws.Range("P15").Value = 0.0300274777826206
$ws.Range("Q15").Value = 45.46971007588667
$ws.Range("R15").Value = 272.81826045532
$ws.Range("S15").Value = 0.009336102011645435
$ws.Range("T15").Value = 0.006629049383122945
$ws.Range("G16").Value = 11.96574466666667
$ws.Range("H16").Value = 35.897234
$ws.Range("I16").Value = 0.2103911152781009
$ws.Range("J16").Value = 0.2207661073338543
$ws.Range("M16").Value = 20.734808
$ws.Range("N16").Value = 62.204424
$ws.Range("O16").Value = 0.242133989171347
$ws.Range("P16").Value = 0.245769325661477
$ws.Range("Q16").Value = 248.1074182403573
$ws.Range("R16").Value = 2232.966764163216
$ws.Range("S16").Value = 0.0509428400284953
$ws.Range("T16").Value = 0.05425753732835061
$ws.Range("G17").Value = 8.018423
$ws.Range("H17").Value = 16.036846
$ws.Range("I17").Value = 0.1409862072722574
$ws.Range("J17").Value = 0.09862576223372788
$ws.Range("M17").Value = 36.48539666666667
$ws.Range("N17").Value = 109.45619
$ws.Range("O17").Value = 0.4260639713374229
$ws.Range("P17").Value = 0.4324607845540777
$ws.Range("Q17").Value = 292.5553437961233
$ws.Range("R17").Value = 1755.33206277674
$ws.Range("S17").Value = 0.06006914337421905
$ws.Range("T17").Value = 0.04265177451284188
$ws.Range("G18").Value = 8.018423
$ws.Range("H18").Value = 16.036846
$ws.Range("I18").Value = 0.1409862072722574
$ws.Range("J18").Value = 0.09862576223372788
$ws.Range("O18").Value = 0.1743777127077069
$ws.Range("P18").Value = 0.1769957741547643
$ws.Range("Q18").Value = 119.7358498336467
$ws.Range("R18").Value = 718.4150990018801
$ws.Range("S18").Value = 0.02458485234747092
$ws.Range("T18").Value = 0.01745634313816238
$ws.Range("G19").Value = 8.018423
$ws.Range("H19").Value = 16.036846
$ws.Range("I19").Value = 0.1409862072722574
$ws.Range("J19").Value = 0.09862576223372788
$ws.Range("M19").Value = 9.680823666666667
$ws.Range("N19").Value = 29.042471
$ws.Range("O19").Value = 0.1130493445068016
$ws.Range("P19").Value = 0.1147466378470605
$ws.Range("Q19").Value = 77.62493914774434
$ws.Range("R19").Value = 465.749634886466
$ws.Range("S19").Value = 0.01593839831662877
$ws.Range("T19").Value = 0.01131697462142387
$ws.Range("G20").Value = 8.018423
$ws.Range("H20").Value = 16.036846
$ws.Range("I20").Value = 0.1409862072722574
$ws.Range("J20").Value = 0.09862576223372788
$ws.Range("M20").Value = 3.79999
$ws.Range("N20").Value = 7.59998
$ws.Range("O20").Value = 0.04437498227672168
$ws.Range("P20").Value = 0.0300274777826206
$ws.Range("Q20").Value = 30.46992721577
$ws.Range("R20").Value = 121.87970886308
$ws.Range("S20").Value = 0.006256260448968634
$ws.Range("T20").Value = 0.002961482884267286
$ws.Range("G21").Value = 8.018423
$ws.Range("H21").Value = 16.036846
$ws.Range("I21").Value = 0.1409862072722574
$ws.Range("J21").Value = 0.09862576223372788
$ws.Range("M21").Value = 20.734808
$ws.Range("N21").Value = 62.204424
$ws.Range("O21").Value = 0.242133989171347
$ws.Range("P21").Value = 0.245769325661477
$ws.Range("Q21").Value = 166.260461367784
$ws.Range("R21").Value = 997.562768206704
$ws.Range("S21").Value = 0.03413755278497007
$ws.Range("T21").Value = 0.02423918707703246
$ws.Range("G22").Value = 7.408770666666666
$ws.Range("H22").Value = 22.226312
$ws.Range("I22").Value = 0.1302668213990815
$ws.Range("J22").Value = 0.1366906536762062
$ws.Range("M22").Value = 36.48539666666667
$ws.Range("N22").Value = 109.45619
$ws.Range("O22").Value = 0.4260639713374229
$ws.Range("P22").Value = 0.4324607845540777
$ws.Range("Q22").Value = 270.3119365856978
$ws.Range("R22").Value = 2432.80742927128
$ws.Range("S22").Value = 0.05550199925879545
$ws.Range("T22").Value = 0.05911334733002185
$ws.Range("G23").Value = 7.408770666666666
$ws.Range("H23").Value = 22.226312
$ws.Range("I23").Value = 0.1302668213990815
$ws.Range("J23").Value = 0.1366906536762062
$ws.Range("O23").Value = 0.1743777127077069
$ws.Range("P23").Value = 0.1769957741547643
$ws.Range("Q23").Value = 110.6321594652622
$ws.Range("R23").Value = 995.6894351873601
$ws.Range("S23").Value = 0.0227156303572752
$ws.Range("T23").Value = 0.0241936680671409
$ws.Range("G24").Value = 7.408770666666666
$ws.Range("H24").Value = 22.226312
$ws.Range("I24").Value = 0.1302668213990815
$ws.Range("J24").Value = 0.1366906536762062
$ws.Range("M24").Value = 9.680823666666667
$ws.Range("N24").Value = 29.042471
$ws.Range("O24").Value = 0.1130493445068016
$ws.Range("P24").Value = 0.1147466378470605
$ws.Range("Q24").Value = 71.72300241077244
$ws.Range("R24").Value = 645.5070216969519
$ws.Range("S24").Value = 0.01472657877015076
$ws.Range("T24").Value = 0.0156847929344616
$ws.Range("G25").Value = 7.408770666666666
$ws.Range("H25").Value = 22.226312
$ws.Range("I25").Value = 0.1302668213990815
$ws.Range("J25").Value = 0.1366906536762062
$ws.Range("M25").Value = 3.79999
$ws.Range("N25").Value = 7.59998
$ws.Range("O25").Value = 0.04437498227672168
$ws.Range("P25").Value = 0.0300274777826206
$ws.Range("Q25").Value = 28.15325444562667
$ws.Range("R25").Value = 168.91952667376
$ws.Range("S25").Value = 0.005780587890829111
$ws.Range("T25").Value = 0.004104475566354168
$ws.Range("G26").Value = 7.408770666666666
$ws.Range("H26").Value = 22.226312
$ws.Range("I26").Value = 0.1302668213990815
$ws.Range("J26").Value = 0.1366906536762062
$ws.Range("M26").Value = 20.734808
$ws.Range("N26").Value = 62.204424
$ws.Range("O26").Value = 0.242133989171347
$ws.Range("P26").Value = 0.245769325661477
$ws.Range("Q26").Value = 153.6194372893653
$ws.Range("R26").Value = 1382.574935604288
$ws.Range("S26").Value = 0.031542025122031
$ws.Range("T26").Value = 0.03359436977822768
